$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the OBI ID list in D18
$ws.Range("D18").Value = "protocol [OBI:0000272]; dose [OBI:0000984]; " + [char]65279 + "extract [OBI:0000423]; organisation [OBI:0000245]; planned process [OBI:0000011]"

# Append a new row (28) for the ADDICTO ontology import
$ws.Range("A28").Value = "ADDICTO"
$ws.Range("B28").Value = "http://purl.obolibrary.org/obo/addicto.owl"
$ws.Range("C28").Value = "entity [BFO:0000001]"
$ws.Range("D28").Value = "substance-induced anxiety disorder [ADDICTO:0001039]"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "all"

# Touch F28 (present-but-empty in the source row) without altering its
# formatting so the cell materializes in the sheet like its neighbours.
$ws.Range("F28").Font.Bold = $false
